$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.538.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.567.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.350'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.84'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.027.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.480.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000143'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.560.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('E21').Value = '  -3.12%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('E24').Value = '  -3.99%  '
$ws.Range('E25').Value = '  -4.30%  '
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.30%  '
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0802'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '456.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.394'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.43'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '158.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.68'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.627'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0532'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('E46').Value = '  -2.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0234'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  -4.06%  '
$ws.Range('E51').Value = '  +3.36%  '
